$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Instrument Properties")
$ws2 = $wb.Worksheets.Item("Curve Properties")

# --- Data changes: rows 134-163 (GBP-USDOIS CrossCurrencySwap block) ---
# Swap the Forecast/Discount curve Left<->Right assignments and
# replace the Convention Left/Right formulas with static convention values.
for ($r = 134; $r -le 163; $r++) {
    $ws1.Cells.Item($r, 4).Value = "na"          # D: Forecast Curve Left
    $ws1.Cells.Item($r, 5).Value = "USDLIBOR3M"  # E: Forecast Curve Right
    $ws1.Cells.Item($r, 6).Value = "GBP-USDOIS"  # F: Discount Curve Left
    $ws1.Cells.Item($r, 7).Value = "USD-USDOIS"  # G: Discount Curve Right
    $ws1.Cells.Item($r, 8).Value = "GBPLIBOR3M"  # H: Convention Left
    $ws1.Cells.Item($r, 9).Value = "USDLIBOR3M"  # I: Convention Right
}

# --- View state changes ---
# Make "Instrument Properties" the active/selected tab (was "Curve Properties").
$ws1.Activate() | Out-Null

$win = $wb.Windows.Item(1)
$win.ScrollRow = 2
$ws1.Range("H21").Select() | Out-Null
